# Auto-generated edit script applying scheduled-runner price/profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1491
$ws.Range("I9").Value = 2103.3333
$ws.Range("J9").Value = 266.33334
$ws.Range("K9").Value = 2103.3333
$ws.Range("L9").Value = 266.33334
$ws.Range("M9").Value = -1934.3333
$ws.Range("N9").Value = -604.33334
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H38").Value = 4718.8335
$ws.Range("I38").Value = 1662.6
$ws.Range("J38").Value = 20000
$ws.Range("K38").Value = 4987.799999999999
$ws.Range("L38").Value = 60000
$ws.Range("M38").Value = -4615.799999999999
$ws.Range("N38").Value = -60744
$ws.Range("H58").Value = 300
$ws.Range("J58").Value = 500
$ws.Range("L58").Value = 1500
$ws.Range("N58").Value = -1800
$ws.Range("H64").Value = 4750
$ws.Range("I64").Value = 4000
$ws.Range("K64").Value = 4000
$ws.Range("M64").Value = -3752
$ws.Range("H67").Value = 4750
$ws.Range("I67").Value = 4000
$ws.Range("K67").Value = 4000
$ws.Range("M67").Value = -3142
$ws.Range("H74").Value = 15366.5
$ws.Range("J74").Value = 4000
$ws.Range("L74").Value = 4000
$ws.Range("N74").Value = -5872
$ws.Range("H76").Value = 2950
$ws.Range("J76").Value = 2900
$ws.Range("L76").Value = 2900
$ws.Range("N76").Value = -3530
$ws.Range("H77").Value = 15366.5
$ws.Range("J77").Value = 4000
$ws.Range("L77").Value = 20000
$ws.Range("N77").Value = -29360
$ws.Range("H79").Value = 2950
$ws.Range("J79").Value = 2900
$ws.Range("L79").Value = 2900
$ws.Range("N79").Value = -5084
$ws.Range("H80").Value = 598.2222
$ws.Range("J80").Value = 548
$ws.Range("L80").Value = 1644
$ws.Range("N80").Value = -3640
$ws.Range("H83").Value = 598.2222
$ws.Range("J83").Value = 548
$ws.Range("L83").Value = 4932
$ws.Range("N83").Value = -14916
$ws.Range("H88").Value = 1466.6666
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 1400
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 1400
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -2212
$ws.Range("H91").Value = 1466.6666
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 1400
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 1400
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -4208
$ws.Range("H94").Value = 1102.5
$ws.Range("I94").Value = 1102.5
$ws.Range("K94").Value = 1102.5
$ws.Range("M94").Value = -651.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 3100
$ws.Range("I9").Value = 3100
$ws.Range("K9").Value = 3100
$ws.Range("M9").Value = -2930
$ws.Range("H20").Value = 3100
$ws.Range("I20").Value = 3100
$ws.Range("K20").Value = 3100
$ws.Range("M20").Value = -2830
$ws.Range("H23").Value = 1999
$ws.Range("J23").Value = 1999
$ws.Range("L23").Value = 1999
$ws.Range("N23").Value = -2517
$ws.Range("H74").Value = 1926.3334
$ws.Range("I74").Value = 1867.25
$ws.Range("K74").Value = 1867.25
$ws.Range("M74").Value = -993.25
$ws.Range("H77").Value = 1926.3334
$ws.Range("I77").Value = 1867.25
$ws.Range("K77").Value = 9336.25
$ws.Range("M77").Value = -4968.25
$ws.Range("H88").Value = 2725.3333
$ws.Range("J88").Value = 2569.2
$ws.Range("L88").Value = 2569.2
$ws.Range("N88").Value = -3381.2
$ws.Range("H91").Value = 2725.3333
$ws.Range("J91").Value = 2569.2
$ws.Range("L91").Value = 2569.2
$ws.Range("N91").Value = -5377.2
$ws.Range("H102").Value = 1698
$ws.Range("I102").Value = 1698
$ws.Range("K102").Value = 1698
$ws.Range("M102").Value = -76

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1323.5
$ws.Range("I86").Value = 1323.5
$ws.Range("K86").Value = 1323.5
$ws.Range("M86").Value = -200.5
$ws.Range("H89").Value = 1323.5
$ws.Range("I89").Value = 1323.5
$ws.Range("K89").Value = 6617.5
$ws.Range("M89").Value = -1001.5
$ws.Range("H107").Value = 931.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 931.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 931.5
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4771.5
$ws.Range("H134").Value = 1981.238
$ws.Range("I134").Value = 1981.238
$ws.Range("K134").Value = 5943.714
$ws.Range("M134").Value = -3408.714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 299
$ws.Range("H16").Value = 3621.1333
$ws.Range("I16").Value = 1919.5
$ws.Range("K16").Value = 1919.5
$ws.Range("M16").Value = -1632.5
$ws.Range("H107").Value = 1500
$ws.Range("I107").Value = 1500
$ws.Range("K107").Value = 1500
$ws.Range("M107").Value = 420
$ws.Range("H113").Value = 3621.1333
$ws.Range("I113").Value = 1919.5
$ws.Range("K113").Value = 1919.5
$ws.Range("M113").Value = 250.5
$ws.Range("H132").Value = 3962.8333
$ws.Range("I132").Value = 4445
$ws.Range("J132").Value = 2998.5
$ws.Range("K132").Value = 13335
$ws.Range("L132").Value = 8995.5
$ws.Range("M132").Value = -10805
$ws.Range("N132").Value = -14055.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 600
$ws.Range("I5").Value = 600
$ws.Range("K5").Value = 1800
$ws.Range("M5").Value = -1688
$ws.Range("H7").Value = 678.8889
$ws.Range("I7").Value = 182
$ws.Range("J7").Value = 1300
$ws.Range("K7").Value = 546
$ws.Range("L7").Value = 3900
$ws.Range("M7").Value = -434
$ws.Range("N7").Value = -4124
$ws.Range("H131").Value = 1949.5
$ws.Range("J131").Value = 3000
$ws.Range("L131").Value = 9000
$ws.Range("N131").Value = -19080
$ws.Range("H135").Value = 600
$ws.Range("I135").Value = 600
$ws.Range("K135").Value = 5400
$ws.Range("M135").Value = -2865

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4734.1665
$ws.Range("I70").Value = 4599.5
$ws.Range("K70").Value = 4599.5
$ws.Range("M70").Value = -4329.5
$ws.Range("H73").Value = 4734.1665
$ws.Range("I73").Value = 4599.5
$ws.Range("K73").Value = 4599.5
$ws.Range("M73").Value = -3663.5
$ws.Range("H74").Value = 49999
$ws.Range("J74").Value = 49999
$ws.Range("L74").Value = 49999
$ws.Range("N74").Value = -51871
$ws.Range("H77").Value = 49999
$ws.Range("J77").Value = 49999
$ws.Range("L77").Value = 149997
$ws.Range("N77").Value = -159357
$ws.Range("H80").Value = 3108.1667
$ws.Range("I80").Value = 2708
$ws.Range("K80").Value = 2708
$ws.Range("M80").Value = -1710
$ws.Range("H83").Value = 3108.1667
$ws.Range("I83").Value = 2708
$ws.Range("K83").Value = 13540
$ws.Range("M83").Value = -8548
$ws.Range("H102").Value = 2689.5557
$ws.Range("I102").Value = 2748
$ws.Range("K102").Value = 2748
$ws.Range("M102").Value = -1126
$ws.Range("H107").Value = 2815.2727
$ws.Range("I107").Value = 1195
$ws.Range("K107").Value = 1195
$ws.Range("M107").Value = 725
$ws.Range("H113").Value = 748.25
$ws.Range("I113").Value = 481
$ws.Range("K113").Value = 481
$ws.Range("M113").Value = 1689

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 33285.715
$ws.Range("J2").Value = 7999.75
$ws.Range("L2").Value = 7999.75
$ws.Range("N2").Value = -8223.75
$ws.Range("H122").Value = 5971.2144
$ws.Range("J122").Value = 7908.364
$ws.Range("L122").Value = 23725.092
$ws.Range("N122").Value = -28625.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2907.3333
$ws.Range("I132").Value = 2907.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8721.999899999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6191.999899999999
$ws.Range("N132").ClearContents()
